# This workbook lists weekly price observations for "Haba" (fava beans) at
# Terminal La Palmera de La Serena. The commit re-dates/re-shuffles the
# per-row observations (date, volume, min/max/avg price, unit price and
# origin) among the existing rows 2-16, without altering the other
# (static/per-market) columns such as A, B, C, E, F, G, H, I, N, Q, R.
#
# Build the new values first (captured from the row each row's data now
# comes from), then write them back in a single pass so that rows don't
# clobber each other while we read them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> values to place into D, J, K, L, M, O, P for that row.
$rowData = @{
    2  = @{ D = 44425; J = 400; K = 11500; L = 12000; M = 11750; O = "Provincia del Elquí"; P = 470 }
    3  = @{ D = 44473; J = 500; K = 8500;  L = 9000;  M = 8750;  O = "Provincia del Elquí"; P = 350 }
    4  = @{ D = 44690; J = 400; K = 17000; L = 18000; M = 17500; O = "Provincia del Elquí"; P = 700 }
    5  = @{ D = 44386; J = 500; K = 11000; L = 12000; M = 11500; O = "Provincia del Elquí"; P = 460 }
    6  = @{ D = 44356; J = 500; K = 13000; L = 14000; M = 13500; O = "Provincia de Limarí"; P = 540 }
    7  = @{ D = 44466; J = 400; K = 9500;  L = 10000; M = 9750;  O = "Provincia del Elquí"; P = 390 }
    8  = @{ D = 44484; J = 400; K = 9000;  L = 10000; M = 9500;  O = "Provincia del Elquí"; P = 380 }
    9  = @{ D = 44446; J = 500; K = 11000; L = 12000; M = 11500; O = "Provincia del Elquí"; P = 460 }
    10 = @{ D = 44370; J = 520; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí"; P = 540 }
    11 = @{ D = 44384; J = 560; K = 11500; L = 12000; M = 11750; O = "Provincia del Elquí"; P = 470 }
    12 = @{ D = 44694; J = 480; K = 17500; L = 18000; M = 17750; O = "Provincia del Elquí"; P = 710 }
    13 = @{ D = 44316; J = 300; K = 16000; L = 17000; M = 16500; O = "Provincia del Elquí"; P = 660 }
    14 = @{ D = 44372; J = 500; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí"; P = 540 }
    15 = @{ D = 44377; J = 520; K = 12500; L = 13000; M = 12750; O = "Provincia del Elquí"; P = 510 }
    16 = @{ D = 44376; J = 400; K = 12000; L = 13000; M = 12500; O = "Provincia del Elquí"; P = 500 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]

    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio $/Kg
}
